$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.Value = "'" + $text
    $c.Style = "Normal"
}

Set-TextValue 'D2' '66.376.66'
Set-TextValue 'E2' '  -1.00%  '
Set-TextValue 'D3' '2.571.40'
Set-TextValue 'E3' '  -2.10%  '
Set-TextValue 'D4' '0.999'
Set-TextValue 'E4' '  -0.11%  '
Set-TextValue 'D5' '582.63'
Set-TextValue 'E5' '  -1.92%  '
Set-TextValue 'D6' '166.21'
Set-TextValue 'E6' '  -0.22%  '
Set-TextValue 'E7' '  -0.05%  '
Set-TextValue 'E8' '  -1.11%  '
Set-TextValue 'D9' '2.570.72'
Set-TextValue 'E9' '  -2.13%  '
Set-TextValue 'E10' '  -2.46%  '
Set-TextValue 'E11' '  +0.15%  '
Set-TextValue 'E12' '  -1.76%  '
Set-TextValue 'E13' '  -1.67%  '
Set-TextValue 'D14' '26.66'
Set-TextValue 'E14' '  -3.51%  '
Set-TextValue 'D15' '3.033.21'
Set-TextValue 'E15' '  -2.63%  '
Set-TextValue 'D16' '0.0000176'
Set-TextValue 'E16' '  -2.65%  '
Set-TextValue 'D17' '66.168.42'
Set-TextValue 'E17' '  -1.00%  '
Set-TextValue 'D18' '2.572.51'
Set-TextValue 'E18' '  -2.14%  '
Set-TextValue 'E19' '  -6.82%  '
Set-TextValue 'D20' '7.72'
Set-TextValue 'E20' '  -4.39%  '
Set-TextValue 'D21' '349.59'
Set-TextValue 'E21' '  -2.22%  '
Set-TextValue 'E22' '  -2.52%  '
Set-TextValue 'D23' '4.59'
Set-TextValue 'E23' '  -1.60%  '
Set-TextValue 'E24' '  +0.01%  '
Set-TextValue 'D25' '1.89'
Set-TextValue 'E25' '  -2.73%  '
Set-TextValue 'D26' '68.91'
Set-TextValue 'E26' '  -1.90%  '
Set-TextValue 'D27' '9.93'
Set-TextValue 'E27' '  -9.19%  '
Set-TextValue 'D28' '2.712.88'
Set-TextValue 'E28' '  -1.93%  '
Set-TextValue 'D29' '0.999'
Set-TextValue 'E29' '  -0.54%  '
Set-TextValue 'D30' '0.0₃0981'
Set-TextValue 'E30' '  -2.79%  '
Set-TextValue 'D31' '527.65'
Set-TextValue 'E31' '  -4.59%  '
Set-TextValue 'D32' '8.13'
Set-TextValue 'E32' '  +2.95%  '
Set-TextValue 'D33' '1.33'
Set-TextValue 'E33' '  -2.60%  '
Set-TextValue 'E34' '  -3.07%  '
Set-TextValue 'E35' '  -3.46%  '
Set-TextValue 'E37' '  -2.94%  '
Set-TextValue 'D38' '156.55'
Set-TextValue 'E38' '  -0.11%  '
Set-TextValue 'D39' '18.73'
Set-TextValue 'E39' '  -1.93%  '
Set-TextValue 'D40' '0.358'
Set-TextValue 'E40' '  -1.93%  '
Set-TextValue 'D41' '18.30'
Set-TextValue 'E41' '  +2.02%  '
Set-TextValue 'E42' '  -1.18%  '
Set-TextValue 'D43' '5.09'
Set-TextValue 'E43' '  -1.27%  '
Set-TextValue 'E44' '  +0.00%  '
Set-TextValue 'D45' '2.43'
Set-TextValue 'E45' '  -0.70%  '
Set-TextValue 'E46' '  -4.33%  '
Set-TextValue 'D47' '148.88'
Set-TextValue 'E47' '  -1.80%  '
Set-TextValue 'E48' '  -3.08%  '
Set-TextValue 'D49' '3.71'
Set-TextValue 'E49' '  -2.24%  '
Set-TextValue 'D50' '1.71'
Set-TextValue 'E50' '  -0.56%  '
Set-TextValue 'D51' '0.0761'
Set-TextValue 'E51' '  -1.32%  '
